# Replace the "TREND IBOVESPA FIA" block (old rows 4-6) with an expanded
# data set covering four funds x four assets (new rows 2-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")

# fund, asset, allowed (kept as literal text "True"/"False", not boolean), dif
$data = @(
    @("TREND DI SIMPLES FUNDOS DE INVESTIMENTO EM RENDA FIXA", "Stocks",  "False"),
    @("TREND DI SIMPLES FUNDOS DE INVESTIMENTO EM RENDA FIXA", "LFT",     "True"),
    @("TREND DI SIMPLES FUNDOS DE INVESTIMENTO EM RENDA FIXA", "Bitcoin", "True"),
    @("TREND DI SIMPLES FUNDOS DE INVESTIMENTO EM RENDA FIXA", "PETR4",   "True"),
    @("EQUITAS SHELTER",                                       "LFT",     "True"),
    @("EQUITAS SHELTER",                                       "Stocks",  "True"),
    @("EQUITAS SHELTER",                                       "Bitcoin", "False"),
    @("EQUITAS SHELTER",                                       "PETR4",   "False"),
    @("DRYS SHELTER PREV",                                     "LFT",     "True"),
    @("DRYS SHELTER PREV",                                     "Stocks",  "True"),
    @("DRYS SHELTER PREV",                                     "Bitcoin", "False"),
    @("DRYS SHELTER PREV",                                     "PETR4",   "False")
)

$rowCount = $data.Count

# Fill column by column (A, then B, then C, then D) so newly-introduced
# shared strings are registered in the same order the source workbook used.
for ($i = 0; $i -lt $rowCount; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $data[$i][0]
}
for ($i = 0; $i -lt $rowCount; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $data[$i][1]
}
for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $i + 2
    # Writing "True"/"False" directly makes the engine coerce the cell to a
    # native boolean (t="b"). The source file stores these as plain text
    # (t="s"), so force text entry with a leading apostrophe and then strip
    # the resulting quote-prefix formatting back to the default style.
    $ws.Cells.Item($row, 3).Value = "'" + $data[$i][2]
    $ws.Cells.Item($row, 3).ClearFormats()
}
for ($i = 0; $i -lt $rowCount; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = 1
}

$ws.Range("B14").Select()
